$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New program entry: "Tracking Bavaria"
# Fill cells in the same order the shared strings were originally appended
# (Links -> short_description/notes -> program_name -> program_identifier)
# so new shared-string indices land in the same order as the source edit.
$ws.Cells.Item(49, 10).Value = "https://ideas.repec.org/p/ces/ifowps/_153.html"
$ws.Cells.Item(49, 6).Value = "In 2000, the Bavarian parliament passed a reform that started tracking (i.e. students attending different levels of secondary education) in 4th instead of 6th grade. The reform only affected the basic and middle track (Hauptschule and Realschule). Piopiunik (2014) finds that the earlier tracking reduces PISA test scores by 13 points on average."
$ws.Cells.Item(49, 2).Value = "Tracking Bavaria"
$ws.Cells.Item(49, 1).Value = "trackingBavaria"

$ws.Cells.Item(49, 3).Value = 2003
$ws.Cells.Item(49, 4).Value = "Education"
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 9).Value = "Piopiunik (2014)"

# Hyperlink for the Links column, then restore the shared Link cell style
# so it matches the style already used by the other hyperlink cells.
$ws.Hyperlinks.Add($ws.Range("J49"), "https://ideas.repec.org/p/ces/ifowps/_153.html")
$ws.Range("J49").Style = "Link"

$ws.Rows.Item(49).RowHeight = 120

# Match the view state shown in the edited workbook
$ws.Activate()
$null = $ws.Range("F49").Select()
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
